$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "h"
$ws.Range("A3").Value = "Firma1"
$ws.Range("A4").Value = "Cp"
$ws.Range("A5").Value = "NIFNIE"
$ws.Range("A6").Value = "Textfield-0"
$ws.Range("A7").Value = "Representante del solicitante"
$ws.Range("A8").Value = "n"
$ws.Range("A9").Value = "m"
$ws.Range("A10").Value = "Textfield-1"
$ws.Range("A11").Value = "Fecha final actuación"
$ws.Range("A13").Value = "AT"
